$d = $word.ActiveDocument

Write-Output ("Before count: " + $d.Paragraphs.Count)

$r = $d.Content
$found = $r.Find.Execute("should be signed by your faculty advisor", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Found: $found"

# Get the paragraph containing this range, then work with paragraph.Range
$para = $r.Paragraphs.First
$pr = $para.Range
$pr.Collapse(1)
$pr.InsertParagraphBefore()

Write-Output ("After count: " + $d.Paragraphs.Count)
